$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# English/Japanese word-list rows were reshuffled (lesson-18 content reorder).
# Row 1 (header) and rows 47-55 are untouched; build the full A2:B113 block
# in its new order and write it back in one shot.
$data = @(
    ('the rest','後|あと'),
    ('air conditioner','エアコン'),
    ('curtain','カーテン'),
    ('cushion','クッション'),
    ('shampoo','シャンプー'),
    ('soy sauce','しょう油|しょうゆ'),
    ('switch','スイッチ'),
    ('soup','スープ'),
    ('skirt','スカート'),
    ('outside','外|そと'),
    ('sofa','ソファ'),
    ('towel','タオル'),
    ('diary','日記|にっき'),
    ('banana','バナナ'),
    ('popcorn','ポップコーン'),
    ('insect','虫|むし'),
    ('rent','家賃|やちん'),
    ('evening','夕方|ゆうがた'),
    ('answering machine','留守番電話|るすばんでんわ'),
    ('refrigerator','冷蔵庫|れいぞうこ'),
    ('candle','ろうそく'),
    ('bright','明るい|あかるい'),
    ('to feel sick','気分が悪い|きぶんがわるい'),
    ('dark','暗い|くらい'),
    ('embarrassing; to feel embarrassed','恥ずかしい|はずかしい'),
    ('(something) opens','開く|あく'),
    ('to apologize','謝る|あやまる'),
    ('to press; to push','押す|おす'),
    ('to drop (something)','落とす|おとす'),
    ('water boils','お湯が沸く|おゆがわく'),
    ('to fall down','転ぶ|ころぶ'),
    ('to break (something)','壊す|こわす'),
    ('to bloom','咲く|さく'),
    ('(something) closes','閉まる|しまる'),
    ('to be saved; to be helped','助かる|たすかる'),
    ('to ask (a favor)','頼む|たのむ'),
    ('(something) turns on','つく'),
    ('to make dirty','汚す|よごす'),
    ('(something) drops','落ちる|おちる'),
    ('to tidy up','片付ける|かたづける'),
    ('to think (about); to consider','考える|かんがえる'),
    ('(something) goes off','消える|きえる'),
    ('(something) breaks','壊れる|こわれる'),
    ('to become dirty','汚れる|よごれる'),
    ('to place an order','注文する|ちゅうもんする'),
    ('right away','今すぐ|いますぐ'),
    ('thanks to...','おかげで'),
    ('See you. (lit., I''m leaving ahead of you.)','お先に失礼します|おさきにしつれいします'),
    ('You must be tired after working so hard. (ritualistic exp.)','お疲れ様（でした）|おつかれさま（でした）'),
    ('short form of ～（ん）でしょう','～（ん）だろう'),
    ('What should I/we do?','どうしよう'),
    ('really','本当に|ほんとうに'),
    ('first of all','まず'),
    ('by (time/date)','～までに'),
    ('purpose','目的|もくてき'),
    ('eye','目|め'),
    ('eye drops','目薬|めぐすり'),
    ('the second','二番目|にばんめ'),
    ('one''s superiors','目上の人|めうえのひと'),
    ('modern','現代的|げんだいてき'),
    ('social','社会的|しゃかいてき'),
    ('target','的|まと'),
    ('physical labor','力仕事|ちからしごと'),
    ('cooperation','協力|きょうりょく'),
    ('endeavor','努力|どりょく'),
    ('sumo wrestler','力士|りきし'),
    ('Western clothes','洋服|ようふく'),
    ('the East','東洋|とうよう'),
    ('Western food','洋食|ようしょく'),
    ('the Atlantic','大西洋|たいせいよう'),
    ('clothes','服|ふく'),
    ('uniform','制服|せいふく'),
    ('Japanese clothes','和服|わふく'),
    ('dining room','食堂|しょくどう'),
    ('public hall','公会堂|こうかいどう'),
    ('dignified; imposing','堂々とした|どうどうとした'),
    ('class','授業|じゅぎょう'),
    ('professor','教授|きょうじゅ'),
    ('to be given','授かる|さずかる'),
    ('class','授業|じゅぎょう'),
    ('occupation','職業|しょくぎょう'),
    ('industry','産業|さんぎょう'),
    ('service industry','サービス業|サービスぎょう'),
    ('exam','試験|しけん'),
    ('game; match','試合|しあい'),
    ('entrance exam','入試|にゅうし'),
    ('to try','試みる|こころみる'),
    ('experiment','実験|じっけん'),
    ('experience','経験|けいけん'),
    ('taking examination','受験|じゅけん'),
    ('to lend','貸す|かす'),
    ('lending','貸し出し|かしだし'),
    ('apartment for rent','賃貸マンション|ちんたいマンション'),
    ('library','図書館|としょかん'),
    ('map','地図|ちず'),
    ('figure','図|ず'),
    ('signal','合図|あいず'),
    ('to attempt','図る|はかる'),
    ('movie theater','映画館|えいがかん'),
    ('Japanese inn','旅館|りょかん'),
    ('embassy','大使館|たいしかん'),
    ('to come to an end','終わる|おわる'),
    ('end','終わり|おわり'),
    ('last stop','終点|しゅうてん'),
    ('the last...','最終～|さいしゅう～'),
    ('homework','宿題|しゅくだい'),
    ('boarding house','下宿|げしゅく'),
    ('lodging','宿泊|しゅくはく'),
    ('inn','宿|やど'),
    ('problem; question','問題|もんだい'),
    ('topic of conversation','話題|わだい'),
    ('title','題|だい')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
